$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update R1's uptime value
$ws.Range("F3").Value = "2:50:00"

# Insert a new row for R3 (Juniper) before the current row 4 (SW1),
# shifting SW1 down to row 5
$ws.Rows("4:4").Insert()

# Populate the new row 4 with the R3 device details
$ws.Range("A4").Value = "R3"
$ws.Range("B4").Value = "R3.automation.local"
$ws.Range("C4").Value = "Juniper"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "12.1R1.9"
$ws.Range("F4").Value = "4:43:56"
$ws.Range("G4").Value = "OK"

# Match style of surrounding data rows (center aligned, like the rest of the table)
$ws.Range("A4:G4").HorizontalAlignment = -4108
$ws.Range("A4:G4").VerticalAlignment = -4108

# Update SW1's uptime value (now on row 5)
$ws.Range("F5").Value = "4:42:00"
